# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect the newly generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 14818
$ws1.Range("F3").Value  = 18334
$ws1.Range("F14").Value = 89
$ws1.Range("F15").Value = 191
$ws1.Range("F17").Value = 1392
$ws1.Range("F22").Value = 7588
$ws1.Range("F26").Value = 1205
$ws1.Range("F28").Value = 5927
$ws1.Range("F29").Value = 91
$ws1.Range("F34").Value = 5254

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 14818
$ws4.Range("F3").Value  = 18334
$ws4.Range("F14").Value = 89
$ws4.Range("F15").Value = 191
$ws4.Range("F17").Value = 1392
$ws4.Range("F23").Value = 7588
$ws4.Range("F27").Value = 1205
$ws4.Range("F30").Value = 5927
$ws4.Range("F31").Value = 91
$ws4.Range("F36").Value = 5254

$wb.Save()
